$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AMSIN": new tenant registration rows for the "176" work-profile runs.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AMSIN")

# Row 15 already existed; align its formatting with the rows above it (style
# "Normal" on every column) and correct the recorded run time in column B.
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2023-04-18"
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").Value = 45034.50914052084
$ws.Range("C15").Value = "176fstrrun"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Style = "Normal"

# Row 16 (new): second work-profile registration run.
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "2023-04-19"
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B16").Value = 45035.68724181713
$ws.Range("C16").Value = "176scndwp"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = 46
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = 46
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Value = 0
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = 0.85
$ws.Range("G16").Style = "Normal"

# Row 17 (new): final work-profile registration run.
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "2023-04-20"
$ws.Range("A17").Style = "Normal"
$ws.Range("B17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B17").Value = 45036.40350466435
$ws.Range("C17").Value = "176fnlwp"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = 46
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = 46
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = 0
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = 0.82
$ws.Range("G17").Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "AMS": new tenant registration rows (beta/live rollout + first work
# profile run).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("AMS")

# Row 13 (new): beta work profile.
$ws2.Range("A13").NumberFormat = "@"
$ws2.Range("A13").Value = "2023-04-20"
$ws2.Range("A13").Style = "Normal"
$ws2.Range("B13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B13").Value = 45036.50684728009
$ws2.Range("C13").Value = "176betawp"
$ws2.Range("C13").Style = "Normal"
$ws2.Range("D13").Value = 46
$ws2.Range("D13").Style = "Normal"
$ws2.Range("E13").Value = 46
$ws2.Range("E13").Style = "Normal"
$ws2.Range("F13").Value = 0
$ws2.Range("F13").Style = "Normal"
$ws2.Range("G13").Value = 0.66
$ws2.Range("G13").Style = "Normal"

# Row 14 (new): live work profile.
$ws2.Range("A14").NumberFormat = "@"
$ws2.Range("A14").Value = "2023-04-27"
$ws2.Range("A14").Style = "Normal"
$ws2.Range("B14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B14").Value = 45043.47313921296
$ws2.Range("C14").Value = "176livewp"
$ws2.Range("C14").Style = "Normal"
$ws2.Range("D14").Value = 46
$ws2.Range("D14").Style = "Normal"
$ws2.Range("E14").Value = 46
$ws2.Range("E14").Style = "Normal"
$ws2.Range("F14").Value = 0
$ws2.Range("F14").Style = "Normal"
$ws2.Range("G14").Value = 0.75
$ws2.Range("G14").Style = "Normal"

# Row 15 (new): first "work profile" tenant run — left in the sheet's
# unformatted default style, matching how brand-new trailing rows land
# before a later formatting pass.
$ws2.Range("A15").NumberFormat = "@"
$ws2.Range("A15").Value = "2023-04-28"
$ws2.Range("A15").Style = "Normal"
$ws2.Range("B15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B15").Value = 45044.70069486451
$ws2.Range("C15").Value = "176fstworkprofile"
$ws2.Range("D15").Value = 46
$ws2.Range("E15").Value = 46
$ws2.Range("F15").Value = 0
$ws2.Range("G15").Value = 1.35
